# Update "want to go" counts (column F) on several sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 869
$ws1.Range("F11").Value = 414
$ws1.Range("F15").Value = 976
$ws1.Range("F17").Value = 406
$ws1.Range("F22").Value = 635
$ws1.Range("F24").Value = 1004

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 340
$ws2.Range("F7").Value = 242
$ws2.Range("F11").Value = 111

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 340
$ws4.Range("F7").Value = 869
$ws4.Range("F15").Value = 414
$ws4.Range("F20").Value = 976
$ws4.Range("F23").Value = 406
$ws4.Range("F26").Value = 242
$ws4.Range("F31").Value = 111
$ws4.Range("F32").Value = 111
$ws4.Range("F34").Value = 635
$ws4.Range("F36").Value = 1004
